# Renamed summary to form
#
# 1. Rename worksheet "Summary" -> "Claim"
# 2. Rename worksheet "Claims"  -> "Submitted"
# 3. Update the Print_Titles defined name that points at the renamed
#    "Claims" sheet so it refers to "Submitted" instead.
# 4. Change the title cell on the Claim sheet from "Claimbot Summary"
#    to "Claimbot" (that cell lives inside a protected range, so the
#    sheet has to be unprotected, edited, then protected again).
# 5. Make "Claim" (the first sheet) the active/selected tab instead of
#    "Insurance".

$wb = $excel.ActiveWorkbook

$wsClaim = $wb.Worksheets.Item(1)
$wsInsurance = $wb.Worksheets.Item(2)
$wsSubmitted = $wb.Worksheets.Item(3)

# --- Rename sheets ------------------------------------------------------
$wsClaim.Name = "Claim"
$wsSubmitted.Name = "Submitted"

# --- Fix up the Print_Titles defined name for the renamed sheet --------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Submitted!Print_Titles") {
        $n.RefersTo = "=Submitted!`$1:`$1"
    }
}

# --- Update the title text on the Claim sheet ---------------------------
$wsClaim.Unprotect()
$wsClaim.Range("A1").Value = "Claimbot"
$wsClaim.Protect()

# --- Make the Claim sheet the selected/active tab -----------------------
$wsClaim.Activate()
